$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.307.20"
$ws.Cells.Item(2, 5).Value = "  -2.29%  "
$ws.Cells.Item(3, 4).Value = "1.872.72"
$ws.Cells.Item(3, 5).Value = "  -1.75%  "
$ws.Cells.Item(4, 5).Value = "  -0.12%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "318.59"
$ws.Cells.Item(5, 5).Value = "  -1.76%  "
$ws.Cells.Item(6, 5).Value = "  -0.09%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4405"
$ws.Cells.Item(7, 5).Value = "  -4.13%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3701"
$ws.Cells.Item(8, 5).Value = "  -3.30%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07518"
$ws.Cells.Item(9, 5).Value = "  -2.56%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.9398"
$ws.Cells.Item(10, 5).Value = "  -4.07%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "21.44"
$ws.Cells.Item(11, 5).Value = "  -2.77%  "
$ws.Cells.Item(12, 4).Value = "1.912.67"
$ws.Cells.Item(12, 5).Value = "  +0.59%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "6.733"
$ws.Cells.Item(13, 5).Value = "  -2.92%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.459"
$ws.Cells.Item(14, 5).Value = "  -3.55%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.06857"
$ws.Cells.Item(15, 5).Value = "  -2.67%  "
$ws.Cells.Item(16, 5).Value = "  -0.14%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "82.23"
$ws.Cells.Item(17, 5).Value = "  -1.83%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000009058"
$ws.Cells.Item(18, 5).Value = "  -4.47%  "
$ws.Cells.Item(19, 5).Value = "  -0.06%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "15.94"
$ws.Cells.Item(20, 5).Value = "  -4.27%  "
$ws.Cells.Item(21, 4).Value = "28.304.12"
$ws.Cells.Item(21, 5).Value = "  -2.33%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.136"
$ws.Cells.Item(22, 5).Value = "  -3.15%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.88"
$ws.Cells.Item(23, 5).Value = "  -0.07%  "
$ws.Cells.Item(24, 4).Value = "2.126.36"
$ws.Cells.Item(24, 5).Value = "  -1.34%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.028"
$ws.Cells.Item(25, 5).Value = "  -3.03%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "154.95"
$ws.Cells.Item(26, 5).Value = "  -1.93%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "18.39"
$ws.Cells.Item(27, 5).Value = "  -3.50%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "5.343"
$ws.Cells.Item(28, 5).Value = "  -5.37%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "113.76"
$ws.Cells.Item(29, 5).Value = "  -3.13%  "
$ws.Cells.Item(30, 5).Value = "  -6.49%  "
$ws.Cells.Item(31, 5).Value = "  -2.47%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.8013"
$ws.Cells.Item(32, 5).Value = "  -7.51%  "
$ws.Cells.Item(33, 5).Value = "  -4.49%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.177"
$ws.Cells.Item(34, 5).Value = "  -5.70%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "2.912"
$ws.Cells.Item(35, 5).Value = "  -1.15%  "
$ws.Cells.Item(36, 5).Value = "  -0.04%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.126"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.05445"
$ws.Cells.Item(38, 5).Value = "  -4.78%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.01959"
$ws.Cells.Item(39, 5).Value = "  -3.96%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.922"
$ws.Cells.Item(40, 5).Value = "  +5.19%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "7.129"
$ws.Cells.Item(41, 5).Value = "  -3.70%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.5265"
$ws.Cells.Item(42, 5).Value = "  -4.28%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.1686"
$ws.Cells.Item(43, 5).Value = "  -3.97%  "
$ws.Cells.Item(44, 5).Value = "  -5.76%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.06756"
$ws.Cells.Item(45, 5).Value = "  -1.17%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.4884"
$ws.Cells.Item(46, 5).Value = "  -5.68%  "
$ws.Cells.Item(47, 2).Value = "RenderToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.992"
$ws.Cells.Item(47, 5).Value = "  -3.46%  "
$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "10.59"
$ws.Cells.Item(48, 5).Value = "  -6.35%  "
$ws.Cells.Item(49, 2).Value = "Quant"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "108.07"
$ws.Cells.Item(49, 5).Value = "  -2.20%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.000002441"
$ws.Cells.Item(50, 5).Value = "  -5.58%  "
$ws.Cells.Item(51, 5).Value = "  -5.11%  "
